$d = $word.ActiveDocument

# The paragraph currently holds three runs: "<id>", "p056v_1", "</id>".
# Collapse them into a single run "<id>p056v_1</id>" (formatting of the
# resulting run follows the first of the matched runs, i.e. the
# Courier New / 7f6000 / 18pt styling used for the "<id>" tag runs).
$d.Content.Find.Execute("<id>p056v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p056v_1</id>", 2) | Out-Null
